# Generate Report for Handoff
# Adds a new handed-off file (9c1c9437-2668-48a8-86fc-59dc4469eee0.md) as a
# new row to the Overview / zh-cn / de-de localization-status worksheets,
# mirroring the existing row for 086322d3-3395-4ce0-b3d8-5ab850d75cfe.md.

$wb = $excel.ActiveWorkbook

$commitHash = "f760515b89f5ecbf0203b8158ffe1475a2b37040"
$newBase    = "9c1c9437-2668-48a8-86fc-59dc4469eee0"
$newMdUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$newBase.md"

$hyperlinkColor = 15570276   # BGR for RGB FF6495ED (matches workbook's "HyperLink" cell style)
$dateFormat     = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 1
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDate($range) {
    $range.NumberFormat = $dateFormat
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

# Seed row 3 from row 2 so untouched columns keep identical values / formats.
$wsOverview.Range("A2:G2").Copy($wsOverview.Range("A3:G3"))

$wsOverview.Range("A3").Value = "$newBase.md"
$wsOverview.Range("B3").Value = "e2e\$newBase.md"
$wsOverview.Range("G3").Value = "2016-09-06 02:43:55"

Style-AsHyperlink $wsOverview.Range("B3")
Style-AsDate $wsOverview.Range("G3")

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdUrl, "", "", "e2e\$newBase.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A2:P2").Copy($wsZh.Range("A3:P3"))

$wsZh.Range("A3").Value = "$newBase.md"
$wsZh.Range("G3").Value = "$newBase.e325eab4d748b93064823db6078a440fce157b98.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-06 02:43:51"

Style-AsHyperlink $wsZh.Range("A3")
Style-AsDate $wsZh.Range("H3")
Style-AsDate $wsZh.Range("K3")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newMdUrl, "", "", "$newBase.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A2:P2").Copy($wsDe.Range("A3:P3"))

$wsDe.Range("A3").Value = "$newBase.md"
$wsDe.Range("G3").Value = "$newBase.e325eab4d748b93064823db6078a440fce157b98.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-06 02:43:55"

Style-AsHyperlink $wsDe.Range("A3")
Style-AsDate $wsDe.Range("H3")
Style-AsDate $wsDe.Range("K3")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newMdUrl, "", "", "$newBase.md") | Out-Null
